$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.337.76'
$ws.Range('E2').Value = '  +0.97%  '
$ws.Range('D3').Value = '1.666.30'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('E4').Value = '  +0.93%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '219.28'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.87%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5340'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.80%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2665'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.49%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06394'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.08%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '20.89'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.56%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07855'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.85%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '4.563'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +1.19%  '
$ws.Range('D13').Value = '1.666.70'
$ws.Range('E13').Value = '  +0.91%  '
$ws.Range('D14').Value = '1.894.99'
$ws.Range('E14').Value = '  +0.86%  '
$ws.Range('E15').Value = '  +0.68%  '
$ws.Range('D16').Value = '0.0₅8189'
$ws.Range('E16').Value = '  -0.27%  '
$ws.Range('E17').Value = '  +0.53%  '
$ws.Range('D18').Value = '26.360.46'
$ws.Range('E18').Value = '  +1.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.683'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.11%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '193.75'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.55%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '10.29'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.06%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.041'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('E24').Value = '  +0.92%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.51'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.67%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.1229'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.67%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.210'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.40%  '
$ws.Range('E28').Value = '  +0.36%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.501'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +4.79%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.05859'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.282'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.77%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.610'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.76%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.282'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.51%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.603'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9696'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +2.44%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.828'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.71%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.422'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.50%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.5830'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.53%  '
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('B40').Value = 'Maker'
$ws.Range('C40').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D40').Value = '1.067.71'
$ws.Range('E40').Value = '  +3.60%  '
$ws.Range('B41').Value = 'TrustWalletToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8630'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.53%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.842'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.63%  '
$ws.Range('E43').Value = '  +0.88%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '104.92'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.76%  '
$ws.Range('D45').Value = '1.806.58'
$ws.Range('E45').Value = '  +0.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '57.81'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.20%  '
$ws.Range('D47').Value = '0.0₈106'
$ws.Range('E47').Value = '  -4.65%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.014'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.00%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.4391'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.976'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.05%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05166'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.44%  '
